$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$modelo = @"
MultiOutputRegressor(estimator=GridSearchCV(cv=5,
                                            estimator=Pipeline(steps=[('model',
                                                                       AdaBoostRegressor())]),
                                            param_grid={'model__learning_rate': [0.1,
                                                                                 0.5,
                                                                                 1.0],
                                                        'model__n_estimators': [50,
                                                                                100,
                                                                                150]},
                                            scoring='neg_mean_squared_error'))
"@

# Trim trailing newline introduced by the here-string
$modelo = $modelo.TrimEnd("`r", "`n")

# New header for column F (copy formatting from A1 so it shares the same style)
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"
$excel.CutCopyMode = $false

# Updated numeric values (B, C, D) for rows 2-10
$ws.Range("B2").Value = 1.361701347923318
$ws.Range("C2").Value = 0.7781647947428578
$ws.Range("D2").Value = 0.978085772356986

$ws.Range("B3").Value = 6.967168367103629
$ws.Range("C3").Value = 0.9004684267492264
$ws.Range("D3").Value = 1.983967610950923

$ws.Range("B4").Value = 4.084931988671991
$ws.Range("C4").Value = 0.7982209121993245
$ws.Range("D4").Value = 1.595290458815578

$ws.Range("B5").Value = 4.102577188453743
$ws.Range("C5").Value = 0.9973052575296558
$ws.Range("D5").Value = 1.717795437638292

$ws.Range("B6").Value = 2.867332255016777
$ws.Range("C6").Value = 0.9687554876831685
$ws.Range("D6").Value = 1.428791619613574

$ws.Range("B7").Value = 2.817726603198308
$ws.Range("C7").Value = 0.9984201887453654
$ws.Range("D7").Value = 1.374402117713929

$ws.Range("B8").Value = 2.050335789537169
$ws.Range("C8").Value = 0.9976514361500552
$ws.Range("D8").Value = 1.148261069783539

$ws.Range("B9").Value = 15.57907473764967
$ws.Range("C9").Value = 0.8137629267999885
$ws.Range("D9").Value = 3.121877164717243

$ws.Range("B10").Value = 1.950958721245652
$ws.Range("C10").Value = 0.994184252479772
$ws.Range("D10").Value = 1.068131567401561

# New column F values (model description) for rows 2-10
$ws.Range("F2").Value = $modelo
$ws.Range("F3").Value = $modelo
$ws.Range("F4").Value = $modelo
$ws.Range("F5").Value = $modelo
$ws.Range("F6").Value = $modelo
$ws.Range("F7").Value = $modelo
$ws.Range("F8").Value = $modelo
$ws.Range("F9").Value = $modelo
$ws.Range("F10").Value = $modelo
